$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp in the title cell ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 03:35"

# --- Guatemala enters the table with fresh data; Tailandia/Senegal/Grecia shift down one row ---
$ws.Range("A77").Value = "Guatemala"
$ws.Range("B77").Value = 3054
$ws.Range("C77").Value = 311
$ws.Range("D77").Value = 244
$ws.Range("E77").Value = 2755
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 55

$ws.Range("A78").Value = "Tailandia"
$ws.Range("B78").Value = 3040
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 2916
$ws.Range("E78").Value = 68
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 56

$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 2976
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 1416
$ws.Range("E79").Value = 1526
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 34

$ws.Range("A80").Value = "Grecia"
$ws.Range("B80").Value = 2876
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 1374
$ws.Range("E80").Value = 1331
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 171

# --- Haiti enters the table with fresh data; Paraguay/Burkina Faso shift down one row ---
$ws.Range("A117").Value = "Haiti"
$ws.Range("B117").Value = 865
$ws.Range("C117").Value = 53
$ws.Range("D117").Value = 22
$ws.Range("E117").Value = 817
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 26

$ws.Range("A118").Value = "Paraguay"
$ws.Range("B118").Value = 850
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 298
$ws.Range("E118").Value = 541
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 11

$ws.Range("A119").Value = "Burkina Faso"
$ws.Range("B119").Value = 814
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 672
$ws.Range("E119").Value = 90
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 52

# --- Nueva Caledonia and Belice swap places (both keep their own stats) ---
$ws.Range("A199").Value = "Nueva Caledonia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Belice"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 16
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 2

# --- Bonaire, San Eustaquio y Saba moves ahead of Sahara Occidental / San Bartolome ---
# (all three rows already shared identical statistics, so only the names move)
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B214").Value = 6
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 6
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("B215").Value = 6
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 6
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "San Bartolome"
$ws.Range("B216").Value = 6
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 6
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
